$wb = $excel.ActiveWorkbook
$grants = $wb.Worksheets.Item(1)

# Insert a brand new row 2 in "grants" holding a valid / acceptable grant
# record, pushing the four existing bad-data rows down one (so the sheet
# now demonstrates one clean row plus the four rows that each contain every
# validation error - hence "all validation errors 4 times").
$grants.Rows.Item(2).Insert()

$grants.Range("A2").Value = "360G-sampletrust-105177/Z/13/Z"
$grants.Range("B2").Value = "Acceptable title"
$grants.Range("C2").Value = "Acceptable description"
$grants.Range("D2").Value = "GBP"
$grants.Range("E2").Value = 1000
$grants.Range("F2").Value = 43617
$grants.Range("F2").NumberFormat = "YYYY\-MM\-DD"

# Widen the Award Date column slightly now that it holds an isoformat date.
$grants.Columns.Item(6).ColumnWidth = 10.15

# "grants" becomes the active tab/selected cell (was "extra sheet" before).
$grants.Range("F3").Select()
